$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Marking" row (row 11): marks for right answer, penalty for wrong answer
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update the "Total" row (row 12): recomputed totals and the summary string
$ws.Range("B12").Value = 55
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "55.0/140"
